# Add Netherlands (and related formatting) to the comparison_rates sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for Netherlands, right after Luxembourg (row 16) ---
# and before New Zealand (old row 17).
$ws.Rows("17:17").Insert()

$ws.Range("A17").Value = "Netherlands"
$ws.Range("B17").Value = 816.8723544279427
$ws.Range("C17").Value = 946.50551065616366
$ws.Range("D17").Value = 1304.3583831012947
$ws.Range("E17").Value = 837.45541638545319
$ws.Range("F17").Value = 515.12122606477283
$ws.Range("G17").Value = 406.35780396617258

# --- Column width tweaks ---
$ws.Columns("B:B").ColumnWidth = 7.666666666666667
$ws.Columns("Q:Q").ColumnWidth = 10.998697916666666
$ws.Columns("R:R").ColumnWidth = 10.666666666666666
$ws.Columns("S:S").ColumnWidth = 10.498697916666666

# --- Apply the "#,##0" number format to the (empty) helper cells in
# columns J:S that trail each data row (leftover chart-range formatting). ---
$ws.Range("J2:S2").NumberFormat = "#,##0"
$ws.Range("J3:S3").NumberFormat = "#,##0"
$ws.Range("J4:S4").NumberFormat = "#,##0"
$ws.Range("J5:S5").NumberFormat = "#,##0"
$ws.Range("J6:S6").NumberFormat = "#,##0"
$ws.Range("J7:S7").NumberFormat = "#,##0"
$ws.Range("J8").NumberFormat = "#,##0"
$ws.Range("L8:S8").NumberFormat = "#,##0"
$ws.Range("J9:N9").NumberFormat = "#,##0"
$ws.Range("P9:S9").NumberFormat = "#,##0"
$ws.Range("J10:S10").NumberFormat = "#,##0"
$ws.Range("J11:S11").NumberFormat = "#,##0"
$ws.Range("J12:S12").NumberFormat = "#,##0"
$ws.Range("J13").NumberFormat = "#,##0"
$ws.Range("Q13:S13").NumberFormat = "#,##0"
$ws.Range("J14:S14").NumberFormat = "#,##0"
$ws.Range("J15:S15").NumberFormat = "#,##0"
$ws.Range("J16").NumberFormat = "#,##0"
$ws.Range("M16").NumberFormat = "#,##0"
$ws.Range("Q16:S16").NumberFormat = "#,##0"
$ws.Range("J17:S17").NumberFormat = "#,##0"
$ws.Range("J18:L18").NumberFormat = "#,##0"
$ws.Range("N18:S18").NumberFormat = "#,##0"
$ws.Range("J19:L19").NumberFormat = "#,##0"
$ws.Range("O19:S19").NumberFormat = "#,##0"
$ws.Range("J20:S20").NumberFormat = "#,##0"
$ws.Range("J21:S21").NumberFormat = "#,##0"
$ws.Range("J22:S22").NumberFormat = "#,##0"
$ws.Range("J23:N23").NumberFormat = "#,##0"
$ws.Range("P23:S23").NumberFormat = "#,##0"
$ws.Range("J24:S24").NumberFormat = "#,##0"
$ws.Range("J25:S25").NumberFormat = "#,##0"
$ws.Range("J26:S26").NumberFormat = "#,##0"

# --- Match the final recorded selection ---
$ws.Range("J32").Select()
